$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Replace the shared "Missing" placeholder text (currently referenced by
#    several C-column cells) with the first new validation-score result, and
#    give the remaining former "Missing" rows their own distinct new scores.
#    We do this BEFORE sorting so the row positions below are the original
#    (pre-sort) ones.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Val Score: 0.5511806337268708, Dice: 0.6550187235030946, F1: 0.5066785952513463, F1_0: 0.7262694768791762, F1_1: 0.5533922882522604, F1_2: 0.4169374062598442, F1_3: 0.43249042623452677"
$ws.Range("C4").Value = "Val Score: 0.5434696009142371, Dice: 0.6444403455583515, F1: 0.5001964246381881, F1_0: 0.7295705033550036, F1_1: 0.5657311055894948, F1_2: 0.40131980139256385, F1_3: 0.42251083235395687"
$ws.Range("C8").Value = "Val Score: 0.5269342885172164, Dice: 0.5924285471543522, F1: 0.49886532052987237, F1_0: 0.7147454741081237, F1_1: 0.5357362183167285, F1_2: 0.4101912815415302, F1_3: 0.43203333773948216 "
$ws.Range("C11").Value = "Val Score: 0.5248356290856423, Dice: 0.6038502582797844, F1: 0.4909722165738671, F1_0: 0.7210279094726519, F1_1: 0.5403606346097436, F1_2: 0.3931981933231123, F1_3: 0.4225923883650231"
$ws.Range("C3").Value = "Val Score: 0.5144322468066552, Dice: 0.6356147390407136, F1: 0.4624968929920588, F1_0: 0.6801651195552217, F1_1: 0.5118687873752198, F1_2: 0.373868221363638, F1_3: 0.3921368319663758"
$ws.Range("C12").Value = "Val Score: 0.5138224758304965, Dice: 0.6061634387730048, F1: 0.4742477774265643, F1_0: 0.7062714425012707, F1_1: 0.48687279436982683, F1_2: 0.38044838280366117, F1_3: 0.4280565530392445"
$ws.Range("C9").Value = "Val Score: 0.5501613593160817, Dice: 0.6525795916374272, F1: 0.5062678311783623, F1_0: 0.728876312405662, F1_1: 0.5447387677327333, F1_2: 0.4135657485837061, F1_3: 0.4395101811662328 "

# ---------------------------------------------------------------------------
# 2. Add the two new training runs (rows appended to the bottom of the
#    existing data, same as the source notebook appending freshly computed
#    rows before re-sorting the sheet).
# ---------------------------------------------------------------------------
$ws.Range("B20").Value = "Add"
$ws.Range("B21").Value = "Add"
$ws.Range("B20").NumberFormat = "0.0000000"
$ws.Range("B21").NumberFormat = "0.0000000"

# ---------------------------------------------------------------------------
# 3. Re-sort the job table by Model, descending, matching the workbook's
#    autoFilter sortState.
# ---------------------------------------------------------------------------
$ws.Range("A2:C24").Sort($ws.Range("B2:B24"), 2, $null, $null, 1)

# ---------------------------------------------------------------------------
# 4. The three now-stale placeholder rows at the bottom of the sorted block
#    are removed (their contents cleared, keeping the sheet's row count and
#    dimension intact).
# ---------------------------------------------------------------------------
$ws.Range("A22:C24").ClearContents()

# ---------------------------------------------------------------------------
# 5. Misc cosmetic bits that round out the commit.
# ---------------------------------------------------------------------------
$ws.Range("D1").Select()
$ws.PageSetup.Orientation = 1
